# Apply the "feedback in French" color highlight.
#
# The commit colors a handful of feedback-list bullets green (RGB
# 00 86 4B) so the translated/ported items stand out. Word stores RGB
# colors as a BGR-packed long (R + G*256 + B*65536), i.e. the same value
# VBA's RGB() function would produce.
$d = $word.ActiveDocument

$green = 0x00 + (0x86 * 256) + (0x4B * 65536)   # => 4949504  (w:color val="00864B")

# Paragraphs whose whole text (run + paragraph mark) must turn green.
$targetTexts = @(
    "Add tagline and adjust image height",
    "? Put href to section below on same page ?",
    "Add “What we do” section in “About Us” page. Picture lady with guy talking on phone and the six bullet points.",
    "Font size of text . Same as the carousel in Main Page (index)",
    "Remove .html"
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]10, [char]7)

    foreach ($target in $targetTexts) {
        if ($text -eq $target) {
            # Paragraph.Range spans the run(s) *and* the paragraph mark,
            # so this colors both <w:r>/<w:rPr> and <w:pPr>/<w:rPr>.
            $para.Range.Font.Color = $green
            Write-Output "Colored paragraph $i : $text"
        }
    }
}
